# Update cryptos list on Sat Jun 22 21:37:19 UTC 2024 with GitHub Actions
#
# Applies the latest price / volume(1h) snapshot to the cryptos sheet,
# including several coins that changed rank position (rows swapped) and
# one newly appeared coin (InjectiveProtocol) at the bottom of the table.

function Set-TextCell($ws, $ref, $val) {
    # Excel's COM `.Value` setter auto-detects numeric-looking strings and
    # converts them to real numbers (losing formatting like trailing
    # zeros, e.g. "1.00" -> 1, "21.17" -> 21.170000000000002). The sheet
    # stores these as plain text, so force text entry (like typing a
    # leading apostrophe in Excel) for anything that looks numeric, then
    # strip the resulting "Text" number-format back to the sheet's default
    # (unstyled) look.
    if ($val -match '^\s*[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?\s*$') {
        $ws.Range($ref).Formula = "'" + $val
        $ws.Range($ref).Style = "Normal"
    } else {
        $ws.Range($ref).Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 'D2' '64.278.71'
Set-TextCell $ws 'E2' '  +0.56%  '
Set-TextCell $ws 'D3' '3.495.32'
Set-TextCell $ws 'E3' '  -0.25%  '
Set-TextCell $ws 'E4' '  +0.04%  '
Set-TextCell $ws 'D5' '588.52'
Set-TextCell $ws 'D6' '134.16'
Set-TextCell $ws 'E6' '  +1.52%  '
Set-TextCell $ws 'E7' '  +0.04%  '
Set-TextCell $ws 'E8' '  +0.59%  '
Set-TextCell $ws 'E9' '  +0.52%  '
Set-TextCell $ws 'E10' '  +2.57%  '
Set-TextCell $ws 'E11' '  +2.75%  '
Set-TextCell $ws 'D12' '4.091.85'
Set-TextCell $ws 'E12' '  -0.25%  '
Set-TextCell $ws 'E13' '  +1.26%  '
Set-TextCell $ws 'E14' '  +1.56%  '
Set-TextCell $ws 'D15' '3.496.36'
Set-TextCell $ws 'E15' '  -0.43%  '
Set-TextCell $ws 'B16' 'WrappedBTC'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D16' '64.366.22'
Set-TextCell $ws 'E16' '  +0.61%  '
Set-TextCell $ws 'B17' 'Avalanche'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D17' '25.67'
Set-TextCell $ws 'E17' '  -6.23%  '
Set-TextCell $ws 'D18' '9.86'
Set-TextCell $ws 'E18' '  +0.71%  '
Set-TextCell $ws 'D19' '5.74'
Set-TextCell $ws 'E19' '  +2.68%  '
Set-TextCell $ws 'E20' '  -2.62%  '
Set-TextCell $ws 'D21' '393.40'
Set-TextCell $ws 'E21' '  +2.89%  '
Set-TextCell $ws 'D22' '0.571'
Set-TextCell $ws 'E22' '  +0.25%  '
Set-TextCell $ws 'D23' '3.636.82'
Set-TextCell $ws 'E23' '  -0.35%  '
Set-TextCell $ws 'D24' '74.58'
Set-TextCell $ws 'E24' '  +1.11%  '
Set-TextCell $ws 'E25' '  +0.04%  '
Set-TextCell $ws 'B26' 'PEPE'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D26' '0.0000116'
Set-TextCell $ws 'E26' '  +1.79%  '
Set-TextCell $ws 'B27' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D27' '1.00'
Set-TextCell $ws 'E27' '  -0.01%  '
Set-TextCell $ws 'B28' 'RenderToken'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 'D28' '7.36'
Set-TextCell $ws 'E28' '  -1.06%  '
Set-TextCell $ws 'B29' 'PancakeSwap'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D29' '2.24'
Set-TextCell $ws 'E29' '  +0.93%  '
Set-TextCell $ws 'B30' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D30' '8.25'
Set-TextCell $ws 'E30' '  -1.37%  '
Set-TextCell $ws 'B31' 'Fetch.AI'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D31' '1.48'
Set-TextCell $ws 'E31' '  -5.67%  '
Set-TextCell $ws 'B32' 'RenzoRestakedETH'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextCell $ws 'D32' '3.519.91'
Set-TextCell $ws 'E32' '  +0.08%  '
Set-TextCell $ws 'B33' 'Kaspa'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D33' '0.152'
Set-TextCell $ws 'E33' '  +5.39%  '
Set-TextCell $ws 'B34' 'USDe'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell $ws 'D34' '1.00'
Set-TextCell $ws 'E34' '  +0.07%  '
Set-TextCell $ws 'B35' 'EthereumClassic'
Set-TextCell $ws 'C35' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D35' '23.44'
Set-TextCell $ws 'E35' '  -0.12%  '
Set-TextCell $ws 'B36' 'NEARProtocol'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D36' '5.13'
Set-TextCell $ws 'E36' '  -3.47%  '
Set-TextCell $ws 'B37' 'Aptos'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D37' '6.88'
Set-TextCell $ws 'E37' '  -0.03%  '
Set-TextCell $ws 'B38' 'ImmutableX'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D38' '1.54'
Set-TextCell $ws 'E38' '  -1.11%  '
Set-TextCell $ws 'B39' 'Monero'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D39' '167.35'
Set-TextCell $ws 'E39' '  +4.45%  '
Set-TextCell $ws 'B40' 'Hedera'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D40' '0.0779'
Set-TextCell $ws 'E40' '  -0.43%  '
Set-TextCell $ws 'B41' 'Mantle'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws 'D41' '0.808'
Set-TextCell $ws 'E41' '  -0.30%  '
Set-TextCell $ws 'B42' 'FirstDigitalUSD'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D42' '1.00'
Set-TextCell $ws 'E42' '  +0.00%  '
Set-TextCell $ws 'B43' 'EnergySwap'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D43' '25.41'
Set-TextCell $ws 'E43' '  -4.77%  '
Set-TextCell $ws 'B44' 'Filecoin'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D44' '4.39'
Set-TextCell $ws 'E44' '  -0.01%  '
Set-TextCell $ws 'B45' 'Stacks'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D45' '1.65'
Set-TextCell $ws 'E45' '  +3.02%  '
Set-TextCell $ws 'B46' 'ONDO'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws 'D46' '1.16'
Set-TextCell $ws 'E46' '  -3.81%  '
Set-TextCell $ws 'B47' 'Maker'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D47' '2.467.04'
Set-TextCell $ws 'E47' '  -0.43%  '
Set-TextCell $ws 'B48' 'Cosmos'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D48' '6.74'
Set-TextCell $ws 'E48' '  -0.48%  '
Set-TextCell $ws 'B49' 'SuiNetwork'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell $ws 'D49' '0.893'
Set-TextCell $ws 'E49' '  -0.93%  '
Set-TextCell $ws 'B50' 'VeChain'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D50' '0.0259'
Set-TextCell $ws 'E50' '  -0.94%  '
Set-TextCell $ws 'B51' 'InjectiveProtocol'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D51' '21.17'
Set-TextCell $ws 'E51' '  -0.10%  '
